$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute(
    "2023-10-28 Saturday", $true, $false, $false, $false, $false,
    $true, 1, $false, "2023-10-29 Sunday", 2)

# Update each division problem in the table, addressed by (row, column)
# so identical/overlapping text values never cross-contaminate each other.
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
$cell.Range.Find.Execute(
    "91÷2=", $true, $false, $false, $false, $false,
    $true, 1, $false, "88÷6=", 2)

$cell = $tbl.Cell(1, 2)
$cell.Range.Find.Execute(
    "20÷3=", $true, $false, $false, $false, $false,
    $true, 1, $false, "24÷9=", 2)

$cell = $tbl.Cell(1, 3)
$cell.Range.Find.Execute(
    "84÷2=", $true, $false, $false, $false, $false,
    $true, 1, $false, "93÷6=", 2)

$cell = $tbl.Cell(1, 4)
$cell.Range.Find.Execute(
    "54÷7=", $true, $false, $false, $false, $false,
    $true, 1, $false, "59÷6=", 2)

$cell = $tbl.Cell(1, 5)
$cell.Range.Find.Execute(
    "10÷3=", $true, $false, $false, $false, $false,
    $true, 1, $false, "49÷6=", 2)

$cell = $tbl.Cell(5, 1)
$cell.Range.Find.Execute(
    "80÷9=", $true, $false, $false, $false, $false,
    $true, 1, $false, "99÷7=", 2)

$cell = $tbl.Cell(5, 2)
$cell.Range.Find.Execute(
    "19÷3=", $true, $false, $false, $false, $false,
    $true, 1, $false, "42÷9=", 2)

$cell = $tbl.Cell(5, 3)
$cell.Range.Find.Execute(
    "32÷9=", $true, $false, $false, $false, $false,
    $true, 1, $false, "51÷2=", 2)

$cell = $tbl.Cell(5, 4)
$cell.Range.Find.Execute(
    "39÷6=", $true, $false, $false, $false, $false,
    $true, 1, $false, "86÷5=", 2)

$cell = $tbl.Cell(5, 5)
$cell.Range.Find.Execute(
    "40÷9=", $true, $false, $false, $false, $false,
    $true, 1, $false, "38÷2=", 2)

$cell = $tbl.Cell(9, 1)
$cell.Range.Find.Execute(
    "59÷8=", $true, $false, $false, $false, $false,
    $true, 1, $false, "62÷6=", 2)

$cell = $tbl.Cell(9, 2)
$cell.Range.Find.Execute(
    "22÷7=", $true, $false, $false, $false, $false,
    $true, 1, $false, "96÷8=", 2)

$cell = $tbl.Cell(9, 3)
$cell.Range.Find.Execute(
    "90÷7=", $true, $false, $false, $false, $false,
    $true, 1, $false, "60÷9=", 2)

$cell = $tbl.Cell(9, 4)
$cell.Range.Find.Execute(
    "31÷4=", $true, $false, $false, $false, $false,
    $true, 1, $false, "30÷8=", 2)

$cell = $tbl.Cell(9, 5)
$cell.Range.Find.Execute(
    "66÷9=", $true, $false, $false, $false, $false,
    $true, 1, $false, "80÷9=", 2)

$cell = $tbl.Cell(13, 1)
$cell.Range.Find.Execute(
    "33÷3=", $true, $false, $false, $false, $false,
    $true, 1, $false, "81÷4=", 2)

$cell = $tbl.Cell(13, 2)
$cell.Range.Find.Execute(
    "70÷4=", $true, $false, $false, $false, $false,
    $true, 1, $false, "80÷6=", 2)

$cell = $tbl.Cell(13, 3)
$cell.Range.Find.Execute(
    "81÷7=", $true, $false, $false, $false, $false,
    $true, 1, $false, "61÷6=", 2)

$cell = $tbl.Cell(13, 4)
$cell.Range.Find.Execute(
    "79÷2=", $true, $false, $false, $false, $false,
    $true, 1, $false, "59÷6=", 2)

$cell = $tbl.Cell(13, 5)
$cell.Range.Find.Execute(
    "18÷3=", $true, $false, $false, $false, $false,
    $true, 1, $false, "46÷4=", 2)

$cell = $tbl.Cell(17, 1)
$cell.Range.Find.Execute(
    "60÷7=", $true, $false, $false, $false, $false,
    $true, 1, $false, "54÷5=", 2)

$cell = $tbl.Cell(17, 2)
$cell.Range.Find.Execute(
    "74÷5=", $true, $false, $false, $false, $false,
    $true, 1, $false, "47÷9=", 2)

$cell = $tbl.Cell(17, 3)
$cell.Range.Find.Execute(
    "62÷5=", $true, $false, $false, $false, $false,
    $true, 1, $false, "30÷4=", 2)

$cell = $tbl.Cell(17, 4)
$cell.Range.Find.Execute(
    "44÷6=", $true, $false, $false, $false, $false,
    $true, 1, $false, "53÷5=", 2)

$cell = $tbl.Cell(17, 5)
$cell.Range.Find.Execute(
    "36÷3=", $true, $false, $false, $false, $false,
    $true, 1, $false, "50÷6=", 2)

